$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3031030526911432
$ws.Range("B3").Value = 0.2646759527437429
$ws.Range("B4").Value = 0.3010672122217473
$ws.Range("B5").Value = 0.1857996930693342
$ws.Range("B6").Value = 0.2666843482755799
$ws.Range("B7").Value = 0.3172781833158986
$ws.Range("B8").Value = 0.2299262201533056
$ws.Range("B9").Value = 0.1585231062812161
$ws.Range("B10").Value = 0.2528836928636566
$ws.Range("B11").Value = 0.2452752211742054
$ws.Range("B12").Value = 0.2330846907205628
$ws.Range("B13").Value = 0.1541328817397854
$ws.Range("B14").Value = 0.1553962895637222
$ws.Range("B15").Value = 0.200170376090257
$ws.Range("B16").Value = 0.1925777228317418
$ws.Range("B17").Value = 0.2769683683347696
$ws.Range("B18").Value = 0.1882054782319281
$ws.Range("B19").Value = 0.2721188087233145
$ws.Range("B20").Value = 0.1911858187425866
$ws.Range("B21").Value = 0.1405487327526666
$ws.Range("B22").Value = 0.1729847822153952
$ws.Range("B23").Value = 0.2356496246225906
$ws.Range("B24").Value = 0.09340496781451638
$ws.Range("B25").Value = 0.1505483795509852
$ws.Range("B26").Value = 0.3400593326103616
$ws.Range("B27").Value = 0.254978200374293
$ws.Range("B28").Value = 0.3759860926245485
$ws.Range("B29").Value = 0.1525970846734771
$ws.Range("B30").Value = 0.1743923826218427
$ws.Range("B31").Value = 0.3441678189982821
$ws.Range("B32").Value = 0.2325964039055812
$ws.Range("B33").Value = 0.2121344956696883
$ws.Range("B34").Value = 0.1657202540739919
$ws.Range("B35").Value = 0.2715570107686192
$ws.Range("B36").Value = 0.2331038502111906
$ws.Range("B37").Value = 0.2015050752556725
$ws.Range("B38").Value = 0.2192844718281782
$ws.Range("B39").Value = 0.2069714739096859
$ws.Range("B40").Value = 0.3383885801018612
$ws.Range("B41").Value = 0.2098195939482268
$ws.Range("B42").Value = 0.2137179949074502
$ws.Range("B43").Value = 0.1795789737747447
$ws.Range("B44").Value = 0.1646274532176297
$ws.Range("B45").Value = 0.2829194960814491
$ws.Range("B46").Value = 0.2141019023010011
$ws.Range("B47").Value = 0.3022534919263937
$ws.Range("B48").Value = 0.2231995161195222
$ws.Range("B49").Value = 0.2879264174706963
$ws.Range("B50").Value = 0.3292303017393478
$ws.Range("B51").Value = 0.3048772031524299
$ws.Range("B52").Value = 0.3199931504041301
$ws.Range("B53").Value = 0.2201985347700067
$ws.Range("B54").Value = 0.190220833396446
$ws.Range("B55").Value = 0.3809247416287085
$ws.Range("B56").Value = 0.3070507106522474
$ws.Range("B57").Value = 0.2051082819470727
$ws.Range("B58").Value = 0.2236368798603683
$ws.Range("B59").Value = 0.2582662286735862
$ws.Range("B60").Value = 0.2627729081303424
$ws.Range("B61").Value = 0.2515903309412535
$ws.Range("B62").Value = 0.2279808317248863
$ws.Range("B63").Value = 0.1880846368812956
$ws.Range("B64").Value = 0.3088656210444864
$ws.Range("B65").Value = 0.4026998536267705
$ws.Range("B66").Value = 0.2882974624343621
$ws.Range("B67").Value = 0.197634552709223
$ws.Range("B68").Value = 0.206237316516559
$ws.Range("B69").Value = 0.2244632121384939
$ws.Range("B70").Value = 0.193780622500697
$ws.Range("B71").Value = 0.188221552986471
$ws.Range("B72").Value = 0.294371349511847
$ws.Range("B73").Value = 0.2608886591279338
$ws.Range("B74").Value = 0.221817571782538
$ws.Range("B75").Value = 0.3541345792767692
$ws.Range("B76").Value = 0.2604765439248585
$ws.Range("B77").Value = 0.2371680202914283
$ws.Range("B78").Value = 0.2062996703527221
